$d = $word.ActiveDocument

# --- Op on P33: delete "Supprimer les espaces vides des formulaires de recherche"
$d.Paragraphs(33).Range.Delete()

# --- Op on P32: delete "Réfléchir sur l’organisation ... (duplicate, to be re-created earlier)"
$d.Paragraphs(32).Range.Delete()

# --- Op on P31: delete "Mettre mes codes css et js dans un fichier séparé ..."
$d.Paragraphs(31).Range.Delete()

# --- Op on P30: delete "Faire pareil si on augmente le nbr de fonctionnalités du site"
$d.Paragraphs(30).Range.Delete()

# --- Op on P29: replace full text -> "Réfléchir sur l’organisation d’un lien..."
$d.Paragraphs(29).Range.Text = "Réfléchir sur l’organisation d’un lien entre les pages (n+1, n-1) entre les reconnaissances, peut-être que les fonction has_next et has_prev peuvent aider à ça ?"

# --- Op on P28: replace full text -> "Créer une page dédiée ... + Présenter fonctionnalités site ..."
$d.Paragraphs(28).Range.Text = "Créer une page dédiée à la présentation du formulaire de recherche, du terrier (photo de l’annexe de mon mémoire) + Présenter fonctionnalités site et penser à mettre un lien vers chacune d’entre elles"

# --- Op on P18: delete "Recherche sur plusieurs pages via woosh ..."
$d.Paragraphs(18).Range.Delete()

# --- Op on P17: delete "Comprendre ce qu’est la recherche plein texte ..."
$d.Paragraphs(17).Range.Delete()

# --- Op on P15: append 4 bits of text at end ("Modifier le readme ...")
$d.Paragraphs(15).Range.InsertAfter(" + Indiquer que le site est moche s’il n’est pas en pleine page + Parler de mes problèmes avec whoosh (genre données rendent difficile navigation vers plusieurs pages + Problème de jsonification des résultats ne permettent pas utilisation dans sa version api)")

# --- Op on P13/P14: P13 full replace, then add new list paragraph after it (taking the old
#     empty P14's place), then delete the old empty paragraph (now shifted to P15)
$d.Paragraphs(13).Range.Text = "Intégrer les données xml de ma charte pour permettre aussi de faire quelques recherches dessus, ça peut s’évaluer -> J’ai un peu bidouillé le truc, pour le moment en vain"
$d.Paragraphs(13).Range.InsertParagraphAfter()
$d.Paragraphs(14).Range.Text = "J’ai un problème avec mon name/331 -> J’aurai voulu que le lieu d’habitation en soit pas « None »"
$d.Paragraphs(15).Range.Delete()

# --- Op on P12: replace full text -> concatenation of 3 parts
$d.Paragraphs(12).Range.Text = "Voir comment on rend le XSLT (site navigable en soi ou devant utiliser l’appli ?) + Penser à la visualisation des images -> Idées de valentin de récupérer via balise image + url construit via les id correspondants, peut-être creuser un peu ça (ou via iiif, mais sans doute plus complexe !) -> Cf devoir de Segolene"

# --- Op on P10: delete "Bon je viens d’intégrer ma charte dans le site ..."
$d.Paragraphs(10).Range.Delete()

# --- Op on P9: append " !" at end ("Penser à nettoyer ...")
$d.Paragraphs(9).Range.InsertAfter(" !")

# --- Op on P8: delete "Tests sont là pour vérifier les éléments clés de mes données"
$d.Paragraphs(8).Range.Delete()

# --- Op on P7: prepend "Tests : " before existing text ("J’ai essayé un truc ...")
$d.Paragraphs(7).Range.InsertBefore("Tests : ")

Write-Output "Done"
